# Refresh the cryptos price list (Price / Volume(1h) columns) on the
# active worksheet, matching the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new Price text, new Volume(1h) text). A $null Price means only
# the Volume(1h) column changed for that row.
$updates = @{
    2  = @("26.050.77", "  +0.92%  ")
    3  = @("1.747.26", "  +0.54%  ")
    4  = @($null, "  -0.02%  ")
    5  = @("232.87", "  +3.24%  ")
    6  = @("1.001", "  +0.08%  ")
    7  = @($null, "  +2.24%  ")
    8  = @("0.2768", "  +2.90%  ")
    9  = @("40.17", "  +2.65%  ")
    10 = @("0.06188", "  +1.96%  ")
    11 = @("1.753.16", "  +0.73%  ")
    12 = @("0.07204", "  +2.94%  ")
    13 = @("15.29", "  +0.75%  ")
    14 = @("0.6399", "  +1.81%  ")
    15 = @("4.581", "  +2.24%  ")
    16 = @("78.14", "  +2.47%  ")
    17 = @("1.000", "  +0.00%  ")
    18 = @("1.001", "  +0.02%  ")
    19 = @("25.982.92", "  +0.58%  ")
    20 = @("11.59", "  +1.62%  ")
    21 = @("0.000006703", "  +2.44%  ")
    22 = @("1.976.98", "  +0.99%  ")
    23 = @("4.326", "  +7.72%  ")
    24 = @("8.784", "  +4.98%  ")
    25 = @("5.200", "  +2.42%  ")
    26 = @("139.14", "  +2.24%  ")
    27 = @("1.521", "  +1.24%  ")
    28 = @("15.28", "  +2.15%  ")
    29 = @("1.811", "  -0.58%  ")
    30 = @("104.32", "  +1.51%  ")
    31 = @("0.08319", "  +0.27%  ")
    32 = @("3.765", "  +4.52%  ")
    33 = @("3.651", "  +9.01%  ")
    34 = @("0.04538", "  +3.22%  ")
    35 = @("2.637", "  +1.12%  ")
    36 = @("0.9989", "  +2.79%  ")
    37 = @("0.6318", "  +6.15%  ")
    38 = @("2.697", "  +0.85%  ")
    39 = @("0.01594", "  +2.19%  ")
    40 = @("1.932", "  +0.44%  ")
    41 = @("1.001", "  +0.15%  ")
    42 = @("98.02", "  -3.70%  ")
    43 = @("0.3901", "  +3.12%  ")
    44 = @("0.7302", "  +0.90%  ")
    45 = @("5.036", "  +4.28%  ")
    46 = @($null, "  +4.24%  ")
    47 = @("6.314", "  +1.04%  ")
    48 = @("0.05345", "  -2.51%  ")
    49 = @("53.92", "  +4.31%  ")
    50 = @("30.53", "  +3.01%  ")
    51 = @("7.681", "  +3.74%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $price = $vals[0]
    $volume = $vals[1]

    if ($null -ne $price) {
        # Several of the refreshed prices read as plain numbers (e.g.
        # "232.87", "1.001"). The source data stores the Price column as
        # text, so force text entry (matching how Excel treats a
        # leading-apostrophe / Text-formatted cell) and then restore the
        # cell's original (default) style so only the value itself
        # changes.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.Style = "Normal"
    }

    $ws.Range("E$row").Value = $volume
}
